$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 28
$ws.Range("A28").Value = 26
$ws.Range("H28").Value = 45535.04839572917
$ws.Range("I28").Value = 116.9734
$ws.Range("J28").Value = 1403.6808
$ws.Range("K28").Value = 0.09
$ws.Range("M28").Value = "completed"
$ws.Range("R28").Value = -27.11559999999998
$ws.Range("W28").Value = "FA1963C27933FB2000"

# Row 29
$ws.Range("A29").Value = 27

# Row 30
$ws.Range("A30").Value = 28

# Row 31
$ws.Range("A31").Value = 29

# Row 32
$ws.Range("A32").Value = 36
$ws.Range("B32").Value = 574
$ws.Range("C32").Value = "ABT"
$ws.Range("D32").Value = 45534.18879386574
$ws.Range("E32").Value = 112.6881
$ws.Range("F32").Value = 1352.2572
$ws.Range("G32").Value = 0.04
$ws.Range("H32").Value = 45535.02311471065
$ws.Range("I32").Value = 113.0516
$ws.Range("J32").Value = 1356.6192
$ws.Range("K32").Value = 0.09
$ws.Range("L32").Value = 12
$ws.Range("M32").Value = "completed"
$ws.Range("R32").Value = 4.231999999999853
$ws.Range("S32").Value = "FA1962B6DF2BFB2000"
$ws.Range("U32").Value = "FA1962BC16BEFB2000"
$ws.Range("V32").Value = "FA1962BC1750FB2000"
$ws.Range("W32").Value = "FA1963BC8A57BB2000"

# Row 33
$ws.Range("A33").Value = 37
$ws.Range("B33").Value = 601
$ws.Range("C33").Value = "PLD"
$ws.Range("D33").Value = 45534.22930445602
$ws.Range("E33").Value = 125.17
$ws.Range("F33").Value = 1376.87
$ws.Range("G33").Value = 0.03
$ws.Range("H33").Value = 45535.01815166666
$ws.Range("I33").Value = 126.05
$ws.Range("J33").Value = 1386.55
$ws.Range("K33").Value = 0.08
$ws.Range("L33").Value = 11
$ws.Range("M33").Value = "completed"
$ws.Range("R33").Value = 9.570000000000064
$ws.Range("S33").Value = "FA1962C4396DBB2000"
$ws.Range("U33").Value = "FA1962C7BC9B3B2000"
$ws.Range("V33").Value = "FA1962C7BD1D04A000"

# Row 34
$ws.Range("A34").Value = 39
$ws.Range("B34").Value = 1
$ws.Range("C34").Value = "CRM"
$ws.Range("D34").Value = 45534.97959670139
$ws.Range("E34").Value = 256.96
$ws.Range("F34").Value = 1284.8
$ws.Range("G34").Value = 0.01
$ws.Range("H34").Value = 45535.04014502315
$ws.Range("I34").Value = 255.9104
$ws.Range("J34").Value = 1279.552
$ws.Range("K34").Value = 0.06999999999999999
$ws.Range("L34").Value = 5
$ws.Range("M34").Value = "completed"
$ws.Range("R34").Value = -5.32799999999982
$ws.Range("S34").Value = "FA1963BB835284A000"
$ws.Range("U34").Value = "FA1963BC88F33B2000"
$ws.Range("V34").Value = "FA1963BC8947FB2000"
$ws.Range("W34").Value = "FA1963C27D0404A000"

# Row 35
$ws.Range("A35").Value = 48
$ws.Range("B35").Value = 10
$ws.Range("C35").Value = "MDT"
$ws.Range("D35").Value = 45534.98202662037
$ws.Range("E35").Value = 88.03740000000001
$ws.Range("F35").Value = 1408.5984
$ws.Range("G35").Value = 0.05
$ws.Range("H35").Value = 45535.04592032408
$ws.Range("I35").Value = 88.3708
$ws.Range("J35").Value = 1413.9328
$ws.Range("K35").Value = 0.09999999999999999
$ws.Range("L35").Value = 16
$ws.Range("M35").Value = "completed"
$ws.Range("R35").Value = 5.18439999999996
$ws.Range("S35").Value = "FA1963BC50597B2000"
$ws.Range("U35").Value = "FA1963C0DA8F7B2000"
$ws.Range("V35").Value = "FA1963C0DAE0C4A000"
$ws.Range("W35").Value = "FA1963C27C4BBB2000"

# Row 36
$ws.Range("A36").Value = 61
$ws.Range("B36").Value = 23
$ws.Range("C36").Value = "ABT"
$ws.Range("D36").Value = 45534.98836896991
$ws.Range("E36").Value = 112.71
$ws.Range("F36").Value = 1352.52
$ws.Range("G36").Value = 0.04
$ws.Range("H36").Value = 45535.04332973379
$ws.Range("I36").Value = 112.81
$ws.Range("J36").Value = 1353.72
$ws.Range("K36").Value = 0.09
$ws.Range("L36").Value = 12
$ws.Range("M36").Value = "completed"
$ws.Range("R36").Value = 1.070000000000045
$ws.Range("S36").Value = "FA1963BE6784BB2000"
$ws.Range("U36").Value = "FA1963C0DBF43B2000"
$ws.Range("V36").Value = "FA1963C0DC46BB2000"
$ws.Range("W36").Value = "FA1963C27DDB3B2000"

# Row 37
$ws.Range("A37").Value = 92
$ws.Range("B37").Value = 13
$ws.Range("C37").Value = "ABT"
$ws.Range("D37").Value = 45535.00843277778
$ws.Range("E37").Value = 112.71
$ws.Range("F37").Value = 1352.52
$ws.Range("G37").Value = 0.04
$ws.Range("H37").Value = 45535.28376490741
$ws.Range("I37").Value = 113.0801
$ws.Range("J37").Value = 1356.9612
$ws.Range("K37").Value = 0.09
$ws.Range("L37").Value = 12
$ws.Range("M37").Value = "completed"
$ws.Range("R37").Value = 4.311199999999981
$ws.Range("S37").Value = "FA1963C5046284A000"
$ws.Range("U37").Value = "FA1963C628D204A000"
$ws.Range("V37").Value = "FA1963C6292644A000"
$ws.Range("W37").Value = "FA1964127CB63B2000"

# Row 38
$ws.Range("A38").Value = 110
$ws.Range("B38").Value = 4
$ws.Range("C38").Value = "ZTS"
$ws.Range("D38").Value = 45535.01617104167
$ws.Range("E38").Value = 182.32
$ws.Range("F38").Value = 1276.24
$ws.Range("G38").Value = 0.02
$ws.Range("H38").Value = 45535.28975931713
$ws.Range("I38").Value = 182.71
$ws.Range("J38").Value = 1278.97
$ws.Range("K38").Value = 0.06999999999999999
$ws.Range("L38").Value = 7
$ws.Range("M38").Value = "completed"
$ws.Range("O38").Value = 0.97
$ws.Range("R38").Value = 2.640000000000018
$ws.Range("S38").Value = "FA1963C7914AC4A000"
$ws.Range("U38").Value = "FA1963CB303EFB2000"
$ws.Range("V38").Value = "FA1963D01F977B2000"
$ws.Range("W38").Value = "FA1964148BC1BB2000"

# Row 39
$ws.Range("A39").Value = 333
$ws.Range("B39").Value = 185
$ws.Range("C39").Value = "SNPS"
$ws.Range("D39").Value = 45535.22942416667
$ws.Range("E39").Value = 515.45
$ws.Range("F39").Value = 2061.8
$ws.Range("G39").Value = 0.01
$ws.Range("H39").Value = 45535.28520369213
$ws.Range("I39").Value = 519.0599999999999
$ws.Range("J39").Value = 2076.24
$ws.Range("K39").Value = 0.07999999999999999
$ws.Range("L39").Value = 4
$ws.Range("M39").Value = "completed"
$ws.Range("O39").Value = 0.97
$ws.Range("R39").Value = 14.3499999999996
$ws.Range("S39").Value = "FA19640DDA66FB2000"
$ws.Range("U39").Value = ""
$ws.Range("V39").Value = ""
$ws.Range("W39").Value = "FA1964135F1BC4A000"
